$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E6").Value = '[''Normal'']'
$ws.Range("D15").Value = '[0, 0, 0, 0, 0, 0, 0]'
$ws.Range("E15").Value = '[]'
$ws.Range("D16").Value = '[1, 0, 0, 0, 1, 0, 0]'
$ws.Range("E16").Value = '[''Normal'', ''RegulationViolation'']'
$ws.Range("D24").Value = '[0, 0, 1, 0, 0, 0, 0]'
$ws.Range("E24").Value = '[''HardwareFault'']'
$ws.Range("D29").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E29").Value = '[''SoftwareFault'']'
$ws.Range("D31").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E31").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D32").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E32").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D38").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E38").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D39").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E39").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D44").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E44").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D45").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E45").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D46").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E46").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D54").Value = '[0, 0, 0, 0, 0, 0, 0]'
$ws.Range("E54").Value = '[]'
$ws.Range("D56").Value = '[0, 0, 1, 0, 0, 0, 0]'
$ws.Range("E56").Value = '[''HardwareFault'']'
$ws.Range("D58").Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Range("E58").Value = '[''Normal'', ''ParamViolation'']'
$ws.Range("D61").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E61").Value = '[''Normal'']'
$ws.Range("D68").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E68").Value = '[''Normal'', ''SurroundingEnvironment'']'
$ws.Range("D73").Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Range("E73").Value = '[''Normal'', ''ParamViolation'']'
$ws.Range("D80").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E80").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D81").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E81").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D84").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E84").Value = '[''Normal'']'
$ws.Range("D88").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E88").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D92").Value = '[1, 0, 1, 0, 0, 0, 1]'
$ws.Range("E92").Value = '[''Normal'', ''HardwareFault'', ''SoftwareFault'']'
$ws.Range("D93").Value = '[1, 0, 1, 0, 0, 0, 1]'
$ws.Range("E93").Value = '[''Normal'', ''HardwareFault'', ''SoftwareFault'']'
$ws.Range("D97").Value = '[0, 0, 1, 0, 0, 0, 0]'
$ws.Range("E97").Value = '[''HardwareFault'']'
$ws.Range("D107").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E107").Value = '[''Normal'']'
$ws.Range("D109").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E109").Value = '[''Normal'', ''SurroundingEnvironment'']'
$ws.Range("D113").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E113").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D116").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E116").Value = '[''Normal'', ''HardwareFault'']'
